$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ST_VAR_01 row (row 2): Stock value 4002 -> 2239
$ws.Range("C2").Value = "1. Menu SP -> Biến thể`n2. Thêm mới`n3. Nhập Stock=2239`n4. Lưu"
$ws.Range("D2").Value = "Stock: 2239"
$ws.Range("F2").Value = "Tìm thấy Stock=2239: true"

# ST_VAR_02 row (row 3): Stock value 4002 -> 2239, 4003 -> 2240
$ws.Range("C3").Value = "1. Tìm Stock=2239`n2. Sửa thành Stock=2240`n3. Lưu"
$ws.Range("D3").Value = "Old Stock: 2239 -> New Stock: 2240"
$ws.Range("F3").Value = "Tìm thấy Stock mới (2240): true"

# ST_VAR_03 row (row 4): Stock value 4003 -> 2240
$ws.Range("C4").Value = "1. Tìm Stock=2240`n2. Xóa`n3. Check biến mất"
$ws.Range("D4").Value = "Target: Stock 2240"
$ws.Range("F4").Value = "Vẫn tìm thấy Stock=2240: false"
